$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1963.5
$ws.Range("I12").Value = 1959.7142
$ws.Range("K12").Value = 1959.7142
$ws.Range("M12").Value = -1789.7142
$ws.Range("H19").Value = 749
$ws.Range("I19").Value = 625
$ws.Range("J19").Value = 831.6667
$ws.Range("K19").Value = 625
$ws.Range("L19").Value = 831.6667
$ws.Range("M19").Value = -450
$ws.Range("N19").Value = -1181.6667
$ws.Range("H105").Value = 40669.332
$ws.Range("J105").Value = 40669.332
$ws.Range("L105").Value = 40669.332
$ws.Range("N105").Value = -47657.332
$ws.Range("H106").Value = 20664
$ws.Range("I106").Value = 2999.5
$ws.Range("K106").Value = 2999.5
$ws.Range("M106").Value = -2368.5
$ws.Range("H132").Value = 1805.4166
$ws.Range("I132").Value = 1623.4117
$ws.Range("K132").Value = 4870.2351
$ws.Range("M132").Value = -2340.2351
$ws.Range("H137").Value = 1793.7
$ws.Range("I137").Value = 1793.7
$ws.Range("K137").Value = 5381.1
$ws.Range("M137").Value = -2831.1
$ws.Range("H139").Value = 148259.67
$ws.Range("J139").Value = 148259.67
$ws.Range("L139").Value = 148259.67
$ws.Range("N139").Value = -158539.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 3764.0417
$ws.Range("I99").Value = 3892.4
$ws.Range("J99").Value = 3550.111
$ws.Range("K99").Value = 3892.4
$ws.Range("L99").Value = 3550.111
$ws.Range("M99").Value = -2394.4
$ws.Range("N99").Value = -6546.111
$ws.Range("H106").Value = 39557
$ws.Range("J106").Value = 39557
$ws.Range("L106").Value = 39557
$ws.Range("N106").Value = -42081
$ws.Range("H107").Value = 1878.2826
$ws.Range("I107").Value = 1463.0968
$ws.Range("K107").Value = 1463.0968
$ws.Range("M107").Value = 456.9032
$ws.Range("H135").Value = 74999
$ws.Range("J135").Value = 74999
$ws.Range("L135").Value = 74999
$ws.Range("N135").Value = -85139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 111111520
$ws.Range("I7").Value = 166667180
$ws.Range("K7").Value = 166667180
$ws.Range("M7").Value = -166667067
$ws.Range("H22").Value = 735.4286
$ws.Range("I22").Value = 799.6667
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 799.6667
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -449.6667
$ws.Range("N22").Value = -1050
$ws.Range("H31").Value = 15714.214
$ws.Range("I31").Value = 4172.3076
$ws.Range("J31").Value = 25717.2
$ws.Range("K31").Value = 4172.3076
$ws.Range("L31").Value = 25717.2
$ws.Range("M31").Value = -3877.3076
$ws.Range("N31").Value = -26307.2
$ws.Range("H34").Value = 15714.214
$ws.Range("I34").Value = 4172.3076
$ws.Range("J34").Value = 25717.2
$ws.Range("K34").Value = 4172.3076
$ws.Range("L34").Value = 25717.2
$ws.Range("M34").Value = -3970.3076
$ws.Range("N34").Value = -26121.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 245
$ws.Range("I26").Value = 205
$ws.Range("J26").Value = 285
$ws.Range("K26").Value = 615
$ws.Range("L26").Value = 855
$ws.Range("M26").Value = -327
$ws.Range("N26").Value = -1431
$ws.Range("H38").Value = 54.75
$ws.Range("I38").Value = 43.166668
$ws.Range("J38").Value = 89.5
$ws.Range("K38").Value = 129.500004
$ws.Range("L38").Value = 268.5
$ws.Range("M38").Value = 217.499996
$ws.Range("N38").Value = -962.5
$ws.Range("H40").Value = 363.33334
$ws.Range("I40").Value = 363.33334
$ws.Range("K40").Value = 1453.33336
$ws.Range("M40").Value = -1384.33336
$ws.Range("H64").Value = 2538.5
$ws.Range("I64").Value = 775
$ws.Range("J64").Value = 2832.4167
$ws.Range("K64").Value = 2325
$ws.Range("L64").Value = 8497.250100000001
$ws.Range("M64").Value = -2055
$ws.Range("N64").Value = -9037.250100000001
$ws.Range("H67").Value = 2538.5
$ws.Range("I67").Value = 775
$ws.Range("J67").Value = 2832.4167
$ws.Range("K67").Value = 2325
$ws.Range("L67").Value = 8497.250100000001
$ws.Range("M67").Value = -1389
$ws.Range("N67").Value = -10369.2501
$ws.Range("H109").Value = 1500
$ws.Range("I109").Value = 1500
$ws.Range("K109").Value = 4500
$ws.Range("M109").Value = -3460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 32786
$ws.Range("J134").Value = 32786
$ws.Range("L134").Value = 98358
$ws.Range("N134").Value = -103428
$ws.Range("H136").Value = 41098.75
$ws.Range("J136").Value = 41098.75
$ws.Range("L136").Value = 123296.25
$ws.Range("N136").Value = -128396.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9205.467000000001
$ws.Range("I7").Value = 13551.077
$ws.Range("K7").Value = 13551.077
$ws.Range("M7").Value = -13439.077
$ws.Range("H68").Value = 1566.0834
$ws.Range("I68").Value = 1454.7778
$ws.Range("K68").Value = 1454.7778
$ws.Range("M68").Value = -705.7778000000001
$ws.Range("H71").Value = 1566.0834
$ws.Range("I71").Value = 1454.7778
$ws.Range("K71").Value = 7273.889
$ws.Range("M71").Value = -3529.889
$ws.Range("H126").Value = 9205.467000000001
$ws.Range("I126").Value = 13551.077
$ws.Range("K126").Value = 40653.231
$ws.Range("M126").Value = -38183.231
$ws.Range("H132").Value = 3376.52
$ws.Range("I132").Value = 2351.4285
$ws.Range("K132").Value = 7054.2855
$ws.Range("M132").Value = -4524.2855
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H139").Value = 109285.57
$ws.Range("J139").Value = 109285.57
$ws.Range("L139").Value = 109285.57
$ws.Range("N139").Value = -119565.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4756.5
$ws.Range("J96").Value = 3564.889
$ws.Range("N96").Value = -6310.889
